$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "non in estrazione" codes on rows 2 and 3
$ws.Range("A2").Value = "39762 (non in estrazione)"
$ws.Range("A3").Value = "39666 (non in estrazione)"

# Append a new row 4 with the same pattern as rows 2/3
$ws.Range("A4").Value = "39742 (non in estrazione)"
$ws.Range("B4").Value = "CAMPO VUOTO"
$ws.Range("C4").Value = "CAMPO VUOTO"
$ws.Range("D4").Value = 0

# Copy the fill/style from B3:C3 onto B4:C4 so the new row matches
$ws.Range("B3:C3").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)
